$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was added to the table. It is inserted at row 48,
# pushing the existing rows 48:185 down to 49:186 (dimension grows from
# A1:R185 to A1:R186).
$ws.Rows.Item(48).Insert()

$ws.Range("A48").Value = 8
$ws.Range("B48").Value = "Terminal La Palmera de La Serena"
$ws.Range("C48").Value = "Coquimbo"
$ws.Range("D48").Value = 44525
$ws.Range("E48").Value = 4
$ws.Range("F48").Value = 100112012
$ws.Range("G48").Value = "Espinaca"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 2000
$ws.Range("K48").Value = 400
$ws.Range("L48").Value = 500
$ws.Range("M48").Value = 450
$ws.Range("N48").Value = "$/atado 300 a 500 gramos"
$ws.Range("O48").Value = "Provincia del Elquí"
$ws.Range("P48").Value = 900
$ws.Range("Q48").Value = 0.5
$ws.Range("R48").Value = "Hortaliza"

# Keep the date column's number format consistent with the rest of column D.
$ws.Range("D48").NumberFormat = $ws.Range("D49").NumberFormat
